$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.166801929473877
$ws.Range("B1").Value = 2.437990665435791
$ws.Range("D1").Value = 2.368453025817871
$ws.Range("E1").Value = 1.234015345573425
